$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 188: correct the timestamp precision and fill in E188/F188
$ws.Range("A188").Value = 45986.43554896991
$ws.Range("E188").Value = 10.9
$ws.Range("F188").Value = 45986
$ws.Range("F188").NumberFormat = "YYYY-MM-DD"

# New data rows 189-197 (same product, price, weight)
$newRows = @(
    @{Row=189; A=45986.48648068287; F=45986},
    @{Row=190; A=45986.49726113426; F=45986},
    @{Row=191; A=45987.42159631944; F=45987},
    @{Row=192; A=45987.42376754629; F=45987},
    @{Row=193; A=45987.42439525463; F=45987},
    @{Row=194; A=45987.42467140046; F=45987},
    @{Row=195; A=45987.4262396875;  F=45987},
    @{Row=196; A=45987.42668586806; F=45987}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("A$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$row").Value = "CREATINA MONOHIDRATO EN POLVO"
    $ws.Range("C$row").Value = "1Kg"
    $ws.Range("D$row").Value = "10,90€"
    $ws.Range("E$row").Value = 10.9
    $ws.Range("F$row").Value = $r.F
    $ws.Range("F$row").NumberFormat = "YYYY-MM-DD"
}

# Row 197: only A-D filled, E and F left blank (inlineStr empty in target)
$ws.Range("A197").Value = 45987.42764274263
$ws.Range("A197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B197").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C197").Value = "1Kg"
$ws.Range("D197").Value = "10,90€"
